$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028821523550997
$ws.Range("D2").Value = 1.031465965228461
$ws.Range("E2").Value = 1.028710766320391
$ws.Range("F2").Value = 1.03807275615048
$ws.Range("I2").Value = 1.032329644113709
$ws.Range("J2").Value = 1.033971715340891
$ws.Range("K2").Value = 1.034274002803758
$ws.Range("L2").Value = 1.031526790337162
$ws.Range("M2").Value = 1.040861830341833
$ws.Range("N2").Value = 1.015238511337861
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03007248251755
$ws.Range("D3").Value = 1.032630246622662
$ws.Range("E3").Value = 1.029782196258577
$ws.Range("F3").Value = 1.039470920725733
$ws.Range("I3").Value = 1.032563765878207
$ws.Range("J3").Value = 1.034861902188711
$ws.Range("K3").Value = 1.03524590045261
$ws.Range("L3").Value = 1.032405502867038
$ws.Range("M3").Value = 1.042068378675414
$ws.Range("N3").Value = 1.015541222587669
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030881240532274
$ws.Range("D4").Value = 1.033383231927748
$ws.Range("E4").Value = 1.030475198248394
$ws.Range("F4").Value = 1.040375058612687
$ws.Range("I4").Value = 1.032713488350325
$ws.Range("J4").Value = 1.035436755945643
$ws.Range("K4").Value = 1.035873832225748
$ws.Range("L4").Value = 1.032973215365554
$ws.Range("M4").Value = 1.042848009353497
$ws.Range("N4").Value = 1.015736506279467
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031221079058205
$ws.Range("D5").Value = 1.033699697793622
$ws.Range("E5").Value = 1.030766470105037
$ws.Range("F5").Value = 1.040755026238178
$ws.Range("I5").Value = 1.032776008573091
$ws.Range("J5").Value = 1.035678149646683
$ws.Range("K5").Value = 1.036137589270729
$ws.Range("L5").Value = 1.03321167452347
$ws.Range("M5").Value = 1.043175509427526
$ws.Range("N5").Value = 1.015818462718887
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031278129959271
$ws.Range("D6").Value = 1.033752828647603
$ws.Range("E6").Value = 1.03081537210877
$ws.Range("F6").Value = 1.040818816854999
$ws.Range("I6").Value = 1.032786481200727
$ws.Range("J6").Value = 1.035718664656288
$ws.Range("K6").Value = 1.036181862081034
$ws.Range("L6").Value = 1.03325170077104
$ws.Range("M6").Value = 1.043230483182036
$ws.Range("N6").Value = 1.015832215324939
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030885782109013
$ws.Range("D7").Value = 1.033387460907429
$ws.Range("E7").Value = 1.03047909049475
$ws.Range("F7").Value = 1.040380136272096
$ws.Range("I7").Value = 1.032714325410337
$ws.Range("J7").Value = 1.035439982537355
$ws.Range("K7").Value = 1.035877357443522
$ws.Range("L7").Value = 1.032976402480221
$ws.Range("M7").Value = 1.042852386429669
$ws.Range("N7").Value = 1.015737601938028
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029244437109992
$ws.Range("D8").Value = 1.031859520271543
$ws.Range("E8").Value = 1.02907292127097
$ws.Range("F8").Value = 1.038545392487591
$ws.Range("I8").Value = 1.032409133522982
$ws.Range("J8").Value = 1.034272798907711
$ws.Range("K8").Value = 1.034602658400263
$ws.Range("L8").Value = 1.031823937304008
$ws.Range("M8").Value = 1.041269816141795
$ws.Range("N8").Value = 1.015340936696992
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026346697999833
$ws.Range("D9").Value = 1.029164053325448
$ws.Range("E9").Value = 1.026592799419222
$ws.Range("F9").Value = 1.03530780639979
$ws.Range("I9").Value = 1.031857762072808
$ws.Range("J9").Value = 1.032207119476185
$ws.Range("K9").Value = 1.032349093291692
$ws.Range("L9").Value = 1.029786373199142
$ws.Range("M9").Value = 1.038472643856605
$ws.Range("N9").Value = 1.014637410542688
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024410978185318
$ws.Range("D10").Value = 1.027364874832375
$ws.Range("E10").Value = 1.024937727986082
$ws.Range("F10").Value = 1.033146110846725
$ws.Range("I10").Value = 1.031481006549958
$ws.Range("J10").Value = 1.030823846642013
$ws.Range("K10").Value = 1.030841612078244
$ws.Range("L10").Value = 1.028423322015123
$ws.Range("M10").Value = 1.036601937115488
$ws.Range("N10").Value = 1.014165294207886
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023571819254811
$ws.Range("D11").Value = 1.026585253146045
$ws.Range("E11").Value = 1.02422064122451
$ws.Range("F11").Value = 1.032209231086892
$ws.Range("I11").Value = 1.031315681577461
$ws.Range("J11").Value = 1.030223384269761
$ws.Range("K11").Value = 1.030187613268001
$ws.Range("L11").Value = 1.027831969257271
$ws.Range("M11").Value = 1.035790444777253
$ws.Range("N11").Value = 1.013960118637935
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023259966937394
$ws.Range("D12").Value = 1.026295578810505
$ws.Range("E12").Value = 1.023954216041298
$ws.Range("F12").Value = 1.031861099252305
$ws.Range("I12").Value = 1.031253943006371
$ws.Range("J12").Value = 1.03000011831362
$ws.Range("K12").Value = 1.02994449840724
$ws.Range("L12").Value = 1.027612140375977
$ws.Range("M12").Value = 1.035488796017534
$ws.Range("N12").Value = 1.013883794391185
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023326867233168
$ws.Range("D13").Value = 1.026357718977826
$ws.Range("E13").Value = 1.024011368232496
$ws.Range("F13").Value = 1.031935780811868
$ws.Range("I13").Value = 1.031267201065088
$ws.Range("J13").Value = 1.030048019957649
$ws.Range("K13").Value = 1.029996656040563
$ws.Range("L13").Value = 1.027659302337209
$ws.Range("M13").Value = 1.035553510925192
$ws.Range("N13").Value = 1.013900171323519
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023546044556202
$ws.Range("D14").Value = 1.026561310391276
$ws.Range("E14").Value = 1.024198619822609
$ws.Range("F14").Value = 1.032180457142994
$ws.Range("I14").Value = 1.031310584975106
$ws.Range("J14").Value = 1.030204933692826
$ws.Range("K14").Value = 1.030167521222223
$ws.Range("L14").Value = 1.02781380170008
$ws.Range("M14").Value = 1.035765514999526
$ws.Range("N14").Value = 1.013953811959843
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023681066761229
$ws.Range("D15").Value = 1.026686738007576
$ws.Range("E15").Value = 1.024313982719155
$ws.Range("F15").Value = 1.032331192532877
$ws.Range("I15").Value = 1.031337271543175
$ws.Range("J15").Value = 1.030301583196787
$ws.Range("K15").Value = 1.030272771559892
$ws.Range("L15").Value = 1.027908970707529
$ws.Range("M15").Value = 1.035896107842197
$ws.Range("N15").Value = 1.013986846740816
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02446664940778
$ws.Range("D16").Value = 1.027416603571409
$ws.Range("E16").Value = 1.024985309389924
$ws.Range("F16").Value = 1.033208270101842
$ws.Range("I16").Value = 1.031491932463127
$ws.Range("J16").Value = 1.030863665598849
$ws.Range("K16").Value = 1.030884989256162
$ws.Range("L16").Value = 1.028462543865986
$ws.Range("M16").Value = 1.036655761967714
$ws.Range("N16").Value = 1.014178895246298
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024959159474248
$ws.Range("D17").Value = 1.027874275320737
$ws.Range("E17").Value = 1.025406297914753
$ws.Range("F17").Value = 1.03375820619025
$ws.Range("I17").Value = 1.031588360930435
$ws.Range("J17").Value = 1.031215842572279
$ws.Range("K17").Value = 1.031268680702217
$ws.Range("L17").Value = 1.028809478034331
$ws.Range("M17").Value = 1.03713187786409
$ws.Range("N17").Value = 1.014299161879734
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.025246338041949
$ws.Range("D18").Value = 1.028141173159695
$ws.Range("E18").Value = 1.025651812099182
$ws.Range("F18").Value = 1.034078892811572
$ws.Range("I18").Value = 1.03164439503791
$ws.Range("J18").Value = 1.031421117275556
$ws.Range("K18").Value = 1.0314923610679
$ws.Range("L18").Value = 1.029011728789543
$ws.Range("M18").Value = 1.037409447260143
$ws.Range("N18").Value = 1.014369239441937
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025344242656056
$ws.Range("D19").Value = 1.028232169301967
$ws.Range("E19").Value = 1.025735519165566
$ws.Range("F19").Value = 1.034188224941075
$ws.Range("I19").Value = 1.031663465460303
$ws.Range("J19").Value = 1.031491086273126
$ws.Range("K19").Value = 1.031568609962374
$ws.Range("L19").Value = 1.029080672473979
$ws.Range("M19").Value = 1.037504067476621
$ws.Range("N19").Value = 1.014393121903442
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.0249063275825
$ws.Range("D20").Value = 1.027825177119083
$ws.Range("E20").Value = 1.025361134122232
$ws.Range("F20").Value = 1.0336992117512
$ws.Range("I20").Value = 1.031578036898614
$ws.Range("J20").Value = 1.031178072248576
$ws.Range("K20").Value = 1.031227526712778
$ws.Range("L20").Value = 1.02877226668591
$ws.Range("M20").Value = 1.037080809746983
$ws.Range("N20").Value = 1.014286265854491
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023481506508397
$ws.Range("D21").Value = 1.02650136025269
$ws.Range("E21").Value = 1.02414348077368
$ws.Range("F21").Value = 1.032108409773592
$ws.Range("I21").Value = 1.031297818598347
$ws.Range("J21").Value = 1.030158732810725
$ws.Range("K21").Value = 1.030117210969021
$ws.Range("L21").Value = 1.027768310315496
$ws.Range("M21").Value = 1.035703091313185
$ws.Range("N21").Value = 1.013938019252696
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022584787335401
$ws.Range("D22").Value = 1.025668512580859
$ws.Range("E22").Value = 1.023377503481812
$ws.Range("F22").Value = 1.031107437533997
$ws.Range("I22").Value = 1.031119727311483
$ws.Range("J22").Value = 1.029516516069676
$ws.Range("K22").Value = 1.029418007306984
$ws.Range("L22").Value = 1.027136074358501
$ws.Range("M22").Value = 1.034835564686297
$ws.Range("N22").Value = 1.013718409152965
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023060239720552
$ws.Range("D23").Value = 1.026110070431628
$ws.Range("E23").Value = 1.023783600419473
$ws.Range("F23").Value = 1.031638146724862
$ws.Range("I23").Value = 1.031214317925303
$ws.Range("J23").Value = 1.029857093119205
$ws.Range("K23").Value = 1.029788774011156
$ws.Range("L23").Value = 1.027471331048162
$ws.Range("M23").Value = 1.035295581708751
$ws.Range("N23").Value = 1.013834890852808
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024930200331206
$ws.Range("D24").Value = 1.027847362650848
$ws.Range("E24").Value = 1.025381541823581
$ws.Range("F24").Value = 1.03372586905129
$ws.Range("I24").Value = 1.031582702536098
$ws.Range("J24").Value = 1.031195139477416
$ws.Range("K24").Value = 1.031246122800966
$ws.Range("L24").Value = 1.028789081232567
$ws.Range("M24").Value = 1.037103885666192
$ws.Range("N24").Value = 1.014292093235437
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027096503888501
$ws.Range("D25").Value = 1.029861272476416
$ws.Range("E25").Value = 1.027234252867081
$ws.Range("F25").Value = 1.036145363591248
$ws.Range("I25").Value = 1.032001918843712
$ws.Range("J25").Value = 1.032742221801606
$ws.Range("K25").Value = 1.032932583206527
$ws.Range("L25").Value = 1.030313947447302
$ws.Range("M25").Value = 1.039196807890675
$ws.Range("N25").Value = 1.014819832249674
